$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing contingency result rows 8-15 ---

# Row 8 (extr1)
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 (extr2)
$ws.Range("C9").Value = 16

# Row 10 (extr3)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11 (extr4)
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12 (extr5)
$ws.Range("C12").Value = 10

# Row 13 (extr6)
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14 (extr7)
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15 (extr8)
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- Append two new rows (line7, line8) at the bottom, matching row 15's format ---

$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A17:E17").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
